$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Reordered "Recorded By" email lists (plain text, no numeric re-parsing
#    risk, so a direct .Value assignment is sufficient).
# ---------------------------------------------------------------------------
$ws.Range("G2").Value  = 'shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G24").Value = 'shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg'

$ws.Range("G18").Value = 'shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G40").Value = 'shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'

$ws.Range("G96").Value  = 'Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Range("G118").Value = 'Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'

$ws.Range("G98").Value  = 'amany.raafat@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Range("G120").Value = 'amany.raafat@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'

$ws.Range("G106").Value = 'wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range("G128").Value = 'wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'

$ws.Range("G134").Value = 'hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg'

$ws.Range("G150").Value = 'wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg'
$ws.Range("G172").Value = 'wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg'

$ws.Range("G156").Value = 'Mohammedeltanany@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg'

# ---------------------------------------------------------------------------
# 2) Class statistics numbers (plain numeric cells).
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 10

$ws.Range("O16").Value = 3
$ws.Range("P16").Value = 2

$ws.Range("O17").Value = 4
$ws.Range("P17").Value = 0

$ws.Range("O21").Value = 4
$ws.Range("P21").Value = 0

# ---------------------------------------------------------------------------
# 3) Percentage values are stored as literal text (e.g. "15.9%"), not real
#    numbers. A plain .Value assignment would be auto-parsed by Excel into a
#    numeric percentage, so force a Text number format first, then restore
#    the original cell style (format only) from an untouched neighbour cell
#    that already carries the same style index.
# ---------------------------------------------------------------------------
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = '15.9%'
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)

$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = '28.4%'
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)

$ws.Range("R16").NumberFormat = "@"
$ws.Range("R16").Value = '13.6%'
$ws.Range("N16").Copy()
$ws.Range("R16").PasteSpecial(-4122)

$ws.Range("S16").NumberFormat = "@"
$ws.Range("S16").Value = '30.0%'
$ws.Range("N16").Copy()
$ws.Range("S16").PasteSpecial(-4122)

$ws.Range("R17").NumberFormat = "@"
$ws.Range("R17").Value = '18.2%'
$ws.Range("N17").Copy()
$ws.Range("R17").PasteSpecial(-4122)

$ws.Range("S17").NumberFormat = "@"
$ws.Range("S17").Value = '35.0%'
$ws.Range("N17").Copy()
$ws.Range("S17").PasteSpecial(-4122)

$ws.Range("R21").NumberFormat = "@"
$ws.Range("R21").Value = '18.2%'
$ws.Range("N21").Copy()
$ws.Range("R21").PasteSpecial(-4122)

$ws.Range("S21").NumberFormat = "@"
$ws.Range("S21").Value = '21.5%'
$ws.Range("N21").Copy()
$ws.Range("S21").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Rows 36 / 54 / 146 changed from "Not Recorded" (pink, style 6) to
#    "Recorded" (green, style 2) now that attendance was taken. Copy the
#    formatting from row 2 (already style 2) onto A:I of each row, then fill
#    in the newly recorded values.
# ---------------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A36:I36").PasteSpecial(-4122)
$ws.Range("G36").Value = 'Amr-Saeed@med.asu.edu.eg'
$ws.Range("H36").Value = '53/217'
$ws.Range("I36").Value = 'Recorded'

$ws.Range("A2:I2").Copy()
$ws.Range("A54:I54").PasteSpecial(-4122)
$ws.Range("G54").Value = 'amany.raafat@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Range("H54").Value = '20/220'
$ws.Range("I54").Value = 'Recorded'

$ws.Range("A2:I2").Copy()
$ws.Range("A146:I146").PasteSpecial(-4122)
$ws.Range("G146").Value = 'Amr-Saeed@med.asu.edu.eg'
$ws.Range("H146").Value = '57/224'
$ws.Range("I146").Value = 'Recorded'
